$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Region" column (C) used the abbreviation "VN" for some Steel rows.
# Update it to match the rest of the sheet, which spells out "Vietnam".
$ws.Range("C36").Value = "Vietnam"
$ws.Range("C39").Value = "Vietnam"
$ws.Range("C40").Value = "Vietnam"
$ws.Range("C41").Value = "Vietnam"
